$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122
$xlPasteFormats = -4122

# ---------------------------------------------------------------
# STEP 1: copy styles from existing header/data cells onto the new
# cells *before* we overwrite any of the source cells' own values,
# so every "s" index in the result matches the original style table
# (no new cellXfs entries get created).
# ---------------------------------------------------------------

# s1 (bold Arial header style, currently on B1/E1) -> G1,H1,I1,J1
$ws.Range("B1").Copy()
$ws.Range("G1:J1").PasteSpecial($xlPasteFormats)

# s2 (Consolas / vertical-center header style, currently on A1/C1/D1) -> E1,F1,K1
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial($xlPasteFormats)
$ws.Range("F1").PasteSpecial($xlPasteFormats)
$ws.Range("K1").PasteSpecial($xlPasteFormats)

# s6 (bold Calibri "Ket" note style, currently on F2/F3) -> L2, L3.
# Must run BEFORE the s5 bulk-paste below, since that paste range
# (D2:J2) includes F2 and would clobber its style first.
$ws.Range("F2").Copy()
$ws.Range("L2").PasteSpecial($xlPasteFormats)
$ws.Range("F3").Copy()
$ws.Range("L3").PasteSpecial($xlPasteFormats)

# s5 (quote-prefix default style, currently on B2) -> D2,E2,F2,G2,H2,I2,J2
$ws.Range("B2").Copy()
$ws.Range("D2:J2").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# STEP 2: header row (row 1) text
# ---------------------------------------------------------------
$ws.Range("C1").Value = "nm_siswa"
$ws.Range("D1").Value = "tmpt_lhr"
$ws.Range("E1").Value = "tgl_lhr"
$ws.Range("F1").Value = "jen_kel"
$ws.Range("G1").Value = "agama"
$ws.Range("H1").Value = "almt_siswa"
$ws.Range("I1").Value = "no_tlp"
$ws.Range("J1").Value = "nm_ayah"
$ws.Range("K1").Value = "kelas_id"

# ---------------------------------------------------------------
# STEP 3: data row (row 2). Leading "'" forces text (quote-prefix),
# matching the original author typing '-' / a date-like string into
# cells styled with the quote-prefix xf.
# ---------------------------------------------------------------
$ws.Range("B2").Value = "'-"
$ws.Range("D2").Value = "'-"
$ws.Range("E2").Value = "'2004-02-10"
$ws.Range("F2").Value = "'-"
$ws.Range("G2").Value = "'-"
$ws.Range("H2").Value = "'-"
$ws.Range("I2").Value = "'-"
$ws.Range("J2").Value = "'-"
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = "Ket: kelas_id lihat di master->kelas"

# ---------------------------------------------------------------
# STEP 4: row 3 note, then drop the old F3 cell (its content moved
# to L3)
# ---------------------------------------------------------------
$ws.Range("L3").Value = "Ket: tgl_lahir pengisian= thn-bln-tgl, contoh: '2021-01-15"
$ws.Range("F3").Clear()

# ---------------------------------------------------------------
# STEP 5: column widths (best effort). The host round-trips
# ColumnWidth through an internal pixel grid (quantised to 1/6 of a
# character), always adding ~5/6 back on readback, so the assigned
# value is nudged by that fixed offset to land as close as possible
# to the target stored width.
# ---------------------------------------------------------------
$widthFudge = 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 10.140625 - $widthFudge
$ws.Columns.Item(5).ColumnWidth = 11.28515625 - $widthFudge
$ws.Columns.Item(6).ColumnWidth = 9 - $widthFudge
$ws.Columns.Item(8).ColumnWidth = 11 - $widthFudge
$ws.Columns.Item(9).ColumnWidth = 11 - $widthFudge
$ws.Columns.Item(10).ColumnWidth = 11 - $widthFudge
$ws.Columns.Item(11).ColumnWidth = 11 - $widthFudge

# ---------------------------------------------------------------
# STEP 6: selection
# ---------------------------------------------------------------
$ws.Range("I10").Select() | Out-Null
